$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Run 50" column (AZ). This shifts the old "Mean" column (BA)
# one position to the left so it becomes the new last column (AZ).
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# Rename the first column header from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Update the "MaxFES" column values (previously generation counts, now fractions)
$colA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
}

# Update the recalculated "Mean" column (now in AZ) values
$colAZ = @(15983553337.13049, 11958397777.3786, 1926011816.208255, 25002.08062499, 1966.27095832, 882.9572044400001, 469.83188608, 249.36217152, 155.65119216, 71.82886056, 44.7303464, 30.66970209, 16.71118063)
for ($i = 0; $i -lt $colAZ.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $colAZ[$i]
}
